$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# 1) "MON Jan 29" / " 10:15:06 PST 2018" were two separate runs; collapse
#    them back down to the single logical sentence (Find/Replace merges
#    the adjacent, identically formatted runs into one run).
# ---------------------------------------------------------------------
$d.Content.Find.Execute("MON Jan 29 10:15:06 PST 2018", $false, $false, $false, $false, $false, $true, 1, $false, "MON Jan 29 10:15:06 PST 2018", 2) | Out-Null

# ---------------------------------------------------------------------
# 2) Append a new purchase-details record (05/02/2018 "SUN Feb 01 ...")
#    right after the very last "Amount Received mode ... - CASH" line.
# ---------------------------------------------------------------------

# Locate the last paragraph whose text is "Amount Received mode<tab><tab>- CASH"
$anchor = $null
for ($i = $d.Paragraphs.Count; $i -ge 1; $i--) {
    $t = $d.Paragraphs($i).Range.Text
    if ($t -like "Amount Received mode*- CASH*") {
        $anchor = $d.Paragraphs($i)
        break
    }
}

function Repeat-Tab([int]$n) {
    $s = ""
    for ($i = 0; $i -lt $n; $i++) {
        $s = $s + "`t"
    }
    return $s
}

function Add-RecordLine([object]$afterPara, [string]$label, [int]$tabCount, [string]$value) {
    $afterPara.Range.InsertParagraphAfter() | Out-Null
    $newIdx = $afterPara.Index + 1
    $newPara = $d.Paragraphs($newIdx)
    $text = $label + (Repeat-Tab $tabCount) + $value
    $newPara.Range.Text = $text
    return $newPara
}

# Blank separator line before the new record.
$anchor.Range.InsertParagraphAfter() | Out-Null
$blank1 = $d.Paragraphs($anchor.Index + 1)

# Timestamp line.
$blank1.Range.InsertParagraphAfter() | Out-Null
$dateLine = $d.Paragraphs($blank1.Index + 1)
$dateLine.Range.Text = "SUN Feb 01 12:07:31 PST 2018"

# Person Name ... - P
$personLine = Add-RecordLine $dateLine "Person Name" 4 "- P"

# Bill number ... - 3006
$billLine = Add-RecordLine $personLine "Bill number" 4 "- 3006"

# Separator dashes.
$billLine.Range.InsertParagraphAfter() | Out-Null
$dashLine = $d.Paragraphs($billLine.Index + 1)
$dashLine.Range.Text = "---------------------------------------------------------------"

# Item Name ... - CARROT EVE
$itemLine = Add-RecordLine $dashLine "Item Name" 4 "- CARROT EVE"

# Number of Pockets ... - 1
$pocketsLine = Add-RecordLine $itemLine "Number of Pockets" 3 "- 1"

# Number of KGs ... - 101
$kgsLine = Add-RecordLine $pocketsLine "Number of KGs" 3 "- 101"

# Rate ... - 12
$rateLine = Add-RecordLine $kgsLine "Rate" 5 "- 12"

# Total Price ... - 1212.0
$totalLine = Add-RecordLine $rateLine "Total Price" 4 "- 1212.0"

# Amount balance ... - 16248.0 (bold)
$totalLine.Range.InsertParagraphAfter() | Out-Null
$balanceLine = $d.Paragraphs($totalLine.Index + 1)
$balanceLine.Range.Text = "Amount balance" + (Repeat-Tab 3) + "- 16248.0"
$balanceLine.Range.Bold = 1

# Two trailing blank lines.
$balanceLine.Range.InsertParagraphAfter() | Out-Null
$blank2 = $d.Paragraphs($balanceLine.Index + 1)

$blank2.Range.InsertParagraphAfter() | Out-Null

Write-Output "done"
